# Append new Lancers listing captured at 2025-10-05 18:23:22 JST.
#
# The tracker keeps the most-recently-seen item pinned at row 2 and inserts
# the single genuinely-new listing right below it (row 3), pushing every
# previously-known listing down by one row. The "fetched at" timestamp in
# column A is refreshed for every row on each run (including the rows that
# merely shifted down), since the whole sheet reflects the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-05 18:23:22"

# 1) Make room for the new listing: push existing rows 3-8 down to rows 4-9.
$ws.Rows.Item(3).Insert()

# 2) Write the brand-new listing into the freshly inserted row 3.
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "初回 【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5407281"
$ws.Range("G3").Value = 158
$ws.Range("H3").Value = "◆自動化,スクレイピング ◇管理"

# 3) Refresh the "fetched at" timestamp on every other data row: row 2
#    (unchanged listing) plus rows 4-9 (listings that shifted down).
$ws.Range("A2").Value = $newTimestamp
for ($r = 4; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 4) The listing that landed at the very bottom (row 9, the old row 8 /
#    "SalesIQ" listing) needs its own hyperlink entry: the row-insert only
#    shifted cell values down, it did not shift the pre-existing hyperlink
#    anchors (which is why F2:F8's hyperlinks still point one row "early").
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5400402")
$ws.Range("F9").Style = $ws.Range("F8").Style
